# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the newly generated site output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows keyed by row number -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 3161
$ws1.Range("F10").Value = 16139
$ws1.Range("F14").Value = 6303
$ws1.Range("F18").Value = 17
$ws1.Range("F21").Value = 15
$ws1.Range("F26").Value = 11
$ws1.Range("F32").Value = 11258
$ws1.Range("F37").Value = 3829

# Sheet "全部类型" - same values, but row numbers shifted by +1 from row 32 onward
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 3161
$ws4.Range("F10").Value = 16139
$ws4.Range("F14").Value = 6303
$ws4.Range("F18").Value = 17
$ws4.Range("F21").Value = 15
$ws4.Range("F26").Value = 11
$ws4.Range("F33").Value = 11258
$ws4.Range("F38").Value = 3829

$wb.Save()
